$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.6811249999999999
$ws.Range("M2").Value = 0.1811433333333334
$ws.Range("N2").Value = 0.5434300000000001
$ws.Range("O2").Value = 0.0111261749556462
$ws.Range("P2").Value = 0.01112617495564619
$ws.Range("Q2").Value = 0.1233812529166667
$ws.Range("R2").Value = 1.11043127625
$ws.Range("S2").Value = 0.0111261749556462
$ws.Range("T2").Value = 0.01112617495564619

$ws.Range("G3").Value = 0.6811249999999999
$ws.Range("O3").Value = 0.8246098959508241
$ws.Range("P3").Value = 0.8246098959508241
$ws.Range("Q3").Value = 9.144328804416665
$ws.Range("R3").Value = 82.29895923974998
$ws.Range("S3").Value = 0.8246098959508241
$ws.Range("T3").Value = 0.8246098959508241

$ws.Range("G4").Value = 0.6811249999999999
$ws.Range("M4").Value = 2.659118666666667
$ws.Range("N4").Value = 7.977356
$ws.Range("O4").Value = 0.1633282272592126
$ws.Range("P4").Value = 0.1633282272592126
$ws.Range("Q4").Value = 1.811192201833333
$ws.Range("R4").Value = 16.3007298165
$ws.Range("S4").Value = 0.1633282272592126
$ws.Range("T4").Value = 0.1633282272592126

$ws.Range("G5").Value = 0.6811249999999999
$ws.Range("M5").Value = 0.015234
$ws.Range("N5").Value = 0.045702
$ws.Range("O5").Value = 0.0009357018343171013
$ws.Range("P5").Value = 0.0009357018343171013
$ws.Range("Q5").Value = 0.01037625825
$ws.Range("R5").Value = 0.09338632424999999
$ws.Range("S5").Value = 0.0009357018343171013
$ws.Range("T5").Value = 0.0009357018343171013
